# Scheduled-runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures across the per-job Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 31271.428
$ws.Range("I10").Value = 31000
$ws.Range("J10").Value = 31380
$ws.Range("K10").Value = 31000
$ws.Range("L10").Value = 31380
$ws.Range("M10").Value = -30707
$ws.Range("N10").Value = -31966

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4921.7
$ws.Range("I61").Value = 3627
$ws.Range("J61").Value = 7942.6665
$ws.Range("K61").Value = 3627
$ws.Range("L61").Value = 7942.6665
$ws.Range("M61").Value = -3415
$ws.Range("N61").Value = -8366.666499999999
$ws.Range("H74").Value = 24851.37
$ws.Range("I74").Value = 26426.896
$ws.Range("J74").Value = 16073.429
$ws.Range("K74").Value = 26426.896
$ws.Range("L74").Value = 16073.429
$ws.Range("M74").Value = -25552.896
$ws.Range("N74").Value = -17821.429
$ws.Range("H77").Value = 24851.37
$ws.Range("I77").Value = 26426.896
$ws.Range("J77").Value = 16073.429
$ws.Range("K77").Value = 132134.48
$ws.Range("L77").Value = 80367.145
$ws.Range("M77").Value = -127766.48
$ws.Range("N77").Value = -89103.145
$ws.Range("H114").Value = 26666.666
$ws.Range("J114").Value = 26666.666
$ws.Range("L114").Value = 26666.666
$ws.Range("N114").Value = -35344.666
$ws.Range("H122").Value = 1848.5883
$ws.Range("I122").Value = 1496.0476
$ws.Range("J122").Value = 2418.077
$ws.Range("K122").Value = 4488.142800000001
$ws.Range("L122").Value = 7254.231000000001
$ws.Range("M122").Value = -2038.142800000001
$ws.Range("N122").Value = -12154.231
$ws.Range("H123").Value = 32426.857
$ws.Range("J123").Value = 32426.857
$ws.Range("L123").Value = 32426.857
$ws.Range("N123").Value = -42226.857
$ws.Range("H132").Value = 2848.2144
$ws.Range("I132").Value = 1931.6471
$ws.Range("J132").Value = 4264.727
$ws.Range("K132").Value = 5794.9413
$ws.Range("L132").Value = 12794.181
$ws.Range("M132").Value = -3264.9413
$ws.Range("N132").Value = -17854.181
$ws.Range("H135").Value = 28620
$ws.Range("J135").Value = 28620
$ws.Range("L135").Value = 28620
$ws.Range("N135").Value = -38760
$ws.Range("H136").Value = 4921.7
$ws.Range("I136").Value = 3627
$ws.Range("J136").Value = 7942.6665
$ws.Range("K136").Value = 10881
$ws.Range("L136").Value = 23827.9995
$ws.Range("M136").Value = -8331
$ws.Range("N136").Value = -28927.9995

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 27999.857
$ws.Range("J40").Value = 27999.857
$ws.Range("L40").Value = 27999.857
$ws.Range("N40").Value = -28529.857
$ws.Range("H96").Value = 16497.375
$ws.Range("I96").Value = 5999.75
$ws.Range("J96").Value = 26995
$ws.Range("K96").Value = 5999.75
$ws.Range("L96").Value = 26995
$ws.Range("M96").Value = -3253.75
$ws.Range("N96").Value = -32487
$ws.Range("H134").Value = 5059.5
$ws.Range("I134").Value = 4842.0586
$ws.Range("J134").Value = 5587.5713
$ws.Range("K134").Value = 14526.1758
$ws.Range("L134").Value = 16762.7139
$ws.Range("M134").Value = -11991.1758
$ws.Range("N134").Value = -21832.7139
$ws.Range("H137").Value = 35814.727
$ws.Range("J137").Value = 35814.727
$ws.Range("L137").Value = 35814.727
$ws.Range("N137").Value = -46014.727

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 27000
$ws.Range("J106").Value = 27000
$ws.Range("L106").Value = 27000
$ws.Range("N106").Value = -29524

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 603.3333
$ws.Range("I6").Value = 46
$ws.Range("J6").Value = 1300
$ws.Range("K6").Value = 138
$ws.Range("L6").Value = 3900
$ws.Range("M6").Value = -25
$ws.Range("N6").Value = -4126
$ws.Range("H7").Value = 16666822
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 18181982
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 54545946
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -54546170
$ws.Range("H17").Value = 1565.3158
$ws.Range("I17").Value = 146.7
$ws.Range("J17").Value = 3141.5557
$ws.Range("K17").Value = 440.1
$ws.Range("L17").Value = 9424.667099999999
$ws.Range("M17").Value = -271.1
$ws.Range("N17").Value = -9762.667099999999
$ws.Range("H122").Value = 1270.2354
$ws.Range("I122").Value = 645
$ws.Range("J122").Value = 1353.6
$ws.Range("K122").Value = 5805
$ws.Range("L122").Value = 12182.4
$ws.Range("M122").Value = -3355
$ws.Range("N122").Value = -17082.4
$ws.Range("H131").Value = 892.88
$ws.Range("I131").Value = 530
$ws.Range("J131").Value = 896.5454999999999
$ws.Range("K131").Value = 1590
$ws.Range("L131").Value = 2689.6365
$ws.Range("M131").Value = 3450
$ws.Range("N131").Value = -12769.6365

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2504551
$ws.Range("I14").Value = 2504551
$ws.Range("K14").Value = 2504551
$ws.Range("M14").Value = -2504383
$ws.Range("H80").Value = 2776.44
$ws.Range("I80").Value = 2700.238
$ws.Range("J80").Value = 3176.5
$ws.Range("K80").Value = 2700.238
$ws.Range("L80").Value = 3176.5
$ws.Range("M80").Value = -1702.238
$ws.Range("N80").Value = -5172.5
$ws.Range("H83").Value = 2776.44
$ws.Range("I83").Value = 2700.238
$ws.Range("J83").Value = 3176.5
$ws.Range("K83").Value = 13501.19
$ws.Range("L83").Value = 15882.5
$ws.Range("M83").Value = -8509.189999999999
$ws.Range("N83").Value = -25866.5
$ws.Range("H123").Value = 12826
$ws.Range("J123").Value = 12826
$ws.Range("L123").Value = 12826
$ws.Range("N123").Value = -17726
$ws.Range("H126").Value = 2224
$ws.Range("I126").Value = 1512.8572
$ws.Range("J126").Value = 2676.5454
$ws.Range("K126").Value = 4538.571599999999
$ws.Range("L126").Value = 8029.6362
$ws.Range("M126").Value = -2068.571599999999
$ws.Range("N126").Value = -12969.6362
$ws.Range("H132").Value = 4323.884
$ws.Range("I132").Value = 2801.8096
$ws.Range("J132").Value = 5776.773
$ws.Range("K132").Value = 8405.4288
$ws.Range("L132").Value = 17330.319
$ws.Range("M132").Value = -5875.4288
$ws.Range("N132").Value = -22390.319

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 4560.6
$ws.Range("J12").Value = 4560.6
$ws.Range("L12").Value = 4560.6
$ws.Range("N12").Value = -4900.6
$ws.Range("H40").Value = 8636.182000000001
$ws.Range("I40").Value = 8599.799999999999
$ws.Range("J40").Value = 9000
$ws.Range("K40").Value = 8599.799999999999
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = -8463.799999999999
$ws.Range("N40").Value = -9272
$ws.Range("H61").Value = 2088.2354
$ws.Range("I61").Value = 1672.7273
$ws.Range("K61").Value = 1672.7273
$ws.Range("M61").Value = -1470.7273
$ws.Range("H96").Value = 26000
$ws.Range("J96").Value = 26000
$ws.Range("L96").Value = 26000
$ws.Range("N96").Value = -31492
$ws.Range("H113").Value = 2088.2354
$ws.Range("I113").Value = 1672.7273
$ws.Range("K113").Value = 1672.7273
$ws.Range("M113").Value = 497.2727
$ws.Range("H122").Value = 2931.4827
$ws.Range("I122").Value = 2389.6428
$ws.Range("J122").Value = 3437.2
$ws.Range("K122").Value = 7168.928400000001
$ws.Range("L122").Value = 10311.6
$ws.Range("M122").Value = -4718.928400000001
$ws.Range("N122").Value = -15211.6
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 30945
$ws.Range("J41").Value = 7338.25
$ws.Range("L41").Value = 7338.25
$ws.Range("N41").Value = -8118.25
$ws.Range("H81").Value = 1806.25
$ws.Range("I81").Value = 1625
$ws.Range("J81").Value = 1987.5
$ws.Range("K81").Value = 3250
$ws.Range("L81").Value = 3975
$ws.Range("M81").Value = -2189
$ws.Range("N81").Value = -6097
$ws.Range("H84").Value = 1806.25
$ws.Range("I84").Value = 1625
$ws.Range("J84").Value = 1987.5
$ws.Range("K84").Value = 16250
$ws.Range("L84").Value = 19875
$ws.Range("M84").Value = -10946
$ws.Range("N84").Value = -30483

# LTW row 133: N133 is removed entirely in the target (no value, cell absent)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N133").ClearContents()

Write-Output "Leve sheets updated."
